$d = $word.ActiveDocument

# Replace each old multiplication expression with its new counterpart.
# All 'old' strings are unique in the document, so a straightforward
# Find/Replace (MatchCase, whole find, no wildcards) for each pair is safe
# and will not cause any chain/duplicate replacements since none of the
# 'new' values coincide with any 'old' value still pending replacement.

$pairs = @(
    @("73×44=", "90×90="),
    @("13×98=", "37×62="),
    @("15×56=", "21×96="),
    @("43×92=", "77×67="),
    @("96×96=", "23×87="),
    @("20×46=", "31×44="),
    @("27×37=", "16×68="),
    @("12×38=", "80×96="),
    @("85×15=", "79×42="),
    @("39×77=", "11×43="),
    @("32×52=", "97×80="),
    @("91×48=", "33×40="),
    @("22×97=", "35×79="),
    @("78×62=", "26×36="),
    @("24×82=", "60×69="),
    @("78×50=", "39×46="),
    @("50×31=", "20×11="),
    @("57×98=", "39×51="),
    @("75×55=", "75×32="),
    @("19×27=", "80×32="),
    @("90×34=", "45×97="),
    @("45×46=", "44×72="),
    @("20×43=", "44×27="),
    @("33×97=", "44×76="),
    @("52×24=", "66×75="),
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}

